# Fix mis-typed cell names: several cells on the "E_I and A_N ratios" sheet
# were labelled "CH_081411_*" when they should have been "CH_081114_*"
# (the day/month digits were transposed). Correct those cell-name values;
# the now-unused "CH_081411_C" / "CH_081411_D" shared strings are dropped
# automatically and the four corrected names are (re)used everywhere they
# belong, including on the "NMDAR" sheet which references the same cells.

$wb = $excel.ActiveWorkbook

$ratios = $wb.Worksheets.Item("E_I and A_N ratios")

$ratios.Range("A16").Value = "CH_081114_A"
$ratios.Range("A17").Value = "CH_081114_A"
$ratios.Range("A18").Value = "CH_081114_B"
$ratios.Range("A19").Value = "CH_081114_B"
$ratios.Range("A20").Value = "CH_081114_C"
$ratios.Range("A21").Value = "CH_081114_D"
$ratios.Range("A22").Value = "CH_081114_D"

# Make "E_I and A_N ratios" the active/selected sheet (it was "NMDAR"
# before), keeping the frozen header pane, and leave the selection on B25.
$ratios.Activate()
$ratios.Range("B25").Select()
